$d = $word.ActiveDocument

# The site footer block ("Ver no Jupiter Salvar em pdf Salvar em docx" /
# "(c) 2020 . Contact: ...") is regenerated by the Jekyll build and was
# dropped from this page. Find that first footer paragraph, then delete it
# together with the blank paragraph right before it and the copyright
# paragraph right after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Ver no Jupiter*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $pBlank = $d.Paragraphs.Item($target - 1)
    $pCopyright = $d.Paragraphs.Item($target + 1)
    $rng = $d.Range($pBlank.Range.Start, $pCopyright.Range.End)
    $rng.Delete()
}
